# Update "data" sheet - Kat's updated L column figures (col L = new enrolment numbers)
# and clear the now-stale M/N derived-split formulas (L4:L13 no longer feed
# the 0.78/0.22 split shown in those columns).
$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

$dataWs.Range("L2").Value = 4216
$dataWs.Range("L3").Value = 4491
$dataWs.Range("L4").Value = 4770
$dataWs.Range("L5").Value = 5059
$dataWs.Range("L6").Value = 5306
$dataWs.Range("L7").Value = 5615
$dataWs.Range("L8").Value = 5899
$dataWs.Range("L9").Value = 6205
$dataWs.Range("L10").Value = 6554
$dataWs.Range("L11").Value = 6864
$dataWs.Range("L12").Value = 7224
$dataWs.Range("L13").Value = 7509

# L14 no longer has a figure for this scenario - clear it along with the
# formulas in M/N that derived from it.
$dataWs.Range("L14").ClearContents()
$dataWs.Range("M4:N14").ClearContents()

# Move the active selection on the data sheet.
[void]$dataWs.Range("M19").Select()

# Update "timepars" sheet column B with Kat's recomputed proportions and
# give those cells a 0.00 number format (matches the new computed-fraction
# look rather than the old rounded literals).
$timeWs = $wb.Worksheets.Item("timepars")

$timeWs.Range("B18").Value = 0.32297447280799113
$timeWs.Range("B18").NumberFormat = "0.00"

$timeWs.Range("B19").Value = 0.35719557195571955
$timeWs.Range("B19").NumberFormat = "0.00"

$timeWs.Range("B20").Value = 0.3888888888888889

$timeWs.Range("B21").Value = 0.49761417859577367
$timeWs.Range("B21").NumberFormat = "0.00"

$timeWs.Range("B22").Value = 0.53368660105980315
$timeWs.Range("B22").NumberFormat = "0.00"

$timeWs.Range("B23").Value = 0.59571788413098237
$timeWs.Range("B23").NumberFormat = "0.00"

# Make "timepars" the active sheet/tab with its own new selection.
[void]$timeWs.Range("F31").Select()
$timeWs.Activate()
